$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 160, shifting existing rows 160:176 down to 161:177.
$ws.Rows("160:160").Insert()

# Populate the newly inserted row 160 with the new record (mirrors the other
# rows for this market/product, only the highlighted fields differ).
$ws.Cells.Item(160, 1).Value = 4
$ws.Cells.Item(160, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(160, 3).Value = "Los Lagos"
$ws.Cells.Item(160, 4).Value = 45154
$ws.Cells.Item(160, 5).Value = 10
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100104
$ws.Cells.Item(160, 8).Value = "Frutos de pepita"
$ws.Cells.Item(160, 9).Value = 100104003
$ws.Cells.Item(160, 10).Value = "Membrillo"
$ws.Cells.Item(160, 11).Value = "Champion"
$ws.Cells.Item(160, 12).Value = "Primera"
$ws.Cells.Item(160, 13).Value = 200
$ws.Cells.Item(160, 14).Value = 15000
$ws.Cells.Item(160, 15).Value = 15000
$ws.Cells.Item(160, 16).Value = 15000
$ws.Cells.Item(160, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(160, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(160, 19).Value = 833
$ws.Cells.Item(160, 20).Value = 18
